$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (pushes old N..P to O..Q),
# matching the "Loan RBI, Variable Instalments" schedule layout change.
$ws.Columns("N").Insert()

# The new (now blank) column N keeps a plain, non bestFit width of 10 chars.
$ws.Columns("N").ColumnWidth = 9.14

# Make "Repayment Schedule" the active sheet/tab and restore the cursor position.
$ws.Activate()
$ws.Range("S8").Select() | Out-Null
